$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ativacao date (B8/C8): "01/01/2018" -> "01/01/2021" ---
# A bare Value assignment of a D/M/Y-shaped string gets auto-parsed into a
# date serial by this host, which would also mint a new cell style. Stage the
# text (quote-prefixed so it is kept literal) in a scratch cell, then
# Copy/PasteSpecial values-only into the targets so they keep their original
# style (s="2"/s="3") and store plain text "01/01/2021".
$scratch = $ws.Range("A1")
$scratch.Value = "'01/01/2021"
$scratch.Copy()
$ws.Range("B8").PasteSpecial(-4163)
$ws.Range("C8").PasteSpecial(-4163)
$scratch.Clear()

# --- Syllabus (B17/C17): short EN text -> full EN syllabus text ---
$syllabusEn = @'
1. Introduction to Operational Research 1.1. Concepts of Operational Research; 1.2. Modeling; 1.3. Structure of Mathematical Models; 1.4. Mathematical techniques in Operational Research; 1.2. Phases of a Study in Operational Research 2. Linear Programming 2.1. Definition 2.2. Formulation of Models 2.3. Graphic Resolution; 3. Simplex method 3.1. Development of the Simplex Method; 3.2. Simplex Method Procedure; 4. Introduction to Graphs and Network Optimization 4.1. Basic Concepts in Graph Theory 4.2. Maximum Flow Problems; 4.3. Minimum Path Problems 5. Case Studies in Linear Programming 5.1. Simple Transport Model 5.2. Model of Designation. 6. Introduction to Queuing Theory 6.1. Queuing Theory Concepts 6.2. Markovian Models
'@
$ws.Range("B17").Value = $syllabusEn
$ws.Range("C17").Value = $syllabusEn

# --- Metodo (B19/C19) ---
$metodoNew = @'
NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.
'@
$ws.Range("B19").Value = $metodoNew
$ws.Range("C19").Value = $metodoNew

# --- Criterio (B20/C20) ---
$criterioNew = @'
NF≥ 5,0.
'@
$ws.Range("B20").Value = $criterioNew
$ws.Range("C20").Value = $criterioNew

# --- Norma de recuperacao (B21/C21) ---
$normaNew = @'
(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.
'@
$ws.Range("B21").Value = $normaNew
$ws.Range("C21").Value = $normaNew

